$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.19105733333333
$ws.Range("H2").Value = 48.573172
$ws.Range("I2").Value = 0.0401918797050022
$ws.Range("J2").Value = 0.0401918797050022
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 2359.049110407361
$ws.Range("R2").Value = 21231.44199366625
$ws.Range("S2").Value = 0.01151874252942787
$ws.Range("T2").Value = 0.01151874252942787

$ws.Range("G3").Value = 16.19105733333333
$ws.Range("H3").Value = 48.573172
$ws.Range("I3").Value = 0.0401918797050022
$ws.Range("J3").Value = 0.0401918797050022
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 2733.04566372562
$ws.Range("R3").Value = 24597.41097353058
$ws.Range("S3").Value = 0.01334488933814037
$ws.Range("T3").Value = 0.01334488933814037

$ws.Range("G4").Value = 16.19105733333333
$ws.Range("H4").Value = 48.573172
$ws.Range("I4").Value = 0.0401918797050022
$ws.Range("J4").Value = 0.0401918797050022
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 2074.497916107534
$ws.Range("R4").Value = 18670.48124496781
$ws.Range("S4").Value = 0.01012933866788006
$ws.Range("T4").Value = 0.01012933866788006

$ws.Range("G5").Value = 16.19105733333333
$ws.Range("H5").Value = 48.573172
$ws.Range("I5").Value = 0.0401918797050022
$ws.Range("J5").Value = 0.0401918797050022
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1064.741400390863
$ws.Range("R5").Value = 9582.672603517763
$ws.Range("S5").Value = 0.005198909169553906
$ws.Range("T5").Value = 0.005198909169553905

$ws.Range("I6").Value = 0.6347354443738135
$ws.Range("J6").Value = 0.6347354443738134
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 37255.58735705774
$ws.Range("R6").Value = 335300.2862135197
$ws.Range("S6").Value = 0.1819112271361119
$ws.Range("T6").Value = 0.1819112271361119

$ws.Range("I7").Value = 0.6347354443738135
$ws.Range("J7").Value = 0.6347354443738134
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.2107508861574761
$ws.Range("T7").Value = 0.2107508861574761

$ws.Range("I8").Value = 0.6347354443738135
$ws.Range("J8").Value = 0.6347354443738134
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 32761.77591836255
$ws.Range("R8").Value = 294855.983265263
$ws.Range("S8").Value = 0.1599688874409501
$ws.Range("T8").Value = 0.15996888744095

$ws.Range("I9").Value = 0.6347354443738135
$ws.Range("J9").Value = 0.6347354443738134
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 16815.06590089088
$ws.Range("R9").Value = 151335.5931080179
$ws.Range("S9").Value = 0.08210444363927545
$ws.Range("T9").Value = 0.08210444363927544

$ws.Range("G10").Value = 57.51647566666667
$ws.Range("H10").Value = 172.549427
$ws.Range("I10").Value = 0.142776053685583
$ws.Range("J10").Value = 0.142776053685583
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 8380.193335235548
$ws.Range("R10").Value = 75421.74001711993
$ws.Range("S10").Value = 0.04091872820686509
$ws.Range("T10").Value = 0.04091872820686508

$ws.Range("G11").Value = 57.51647566666667
$ws.Range("H11").Value = 172.549427
$ws.Range("I11").Value = 0.142776053685583
$ws.Range("J11").Value = 0.142776053685583
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 9708.763990967904
$ws.Range("R11").Value = 87378.87591871113
$ws.Range("S11").Value = 0.04740586035177052
$ws.Range("T11").Value = 0.04740586035177051

$ws.Range("G12").Value = 57.51647566666667
$ws.Range("H12").Value = 172.549427
$ws.Range("I12").Value = 0.142776053685583
$ws.Range("J12").Value = 0.142776053685583
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 7369.364857148905
$ws.Range("R12").Value = 66324.28371434013
$ws.Range("S12").Value = 0.03598306454088788
$ws.Range("T12").Value = 0.03598306454088787

$ws.Range("G13").Value = 57.51647566666667
$ws.Range("H13").Value = 172.549427
$ws.Range("I13").Value = 0.142776053685583
$ws.Range("J13").Value = 0.142776053685583
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 3782.345500117244
$ws.Range("R13").Value = 34041.1095010552
$ws.Range("S13").Value = 0.01846840058605957
$ws.Range("T13").Value = 0.01846840058605957

$ws.Range("G14").Value = 73.43709933333334
$ws.Range("H14").Value = 220.311298
$ws.Range("I14").Value = 0.1822966222356014
$ws.Range("J14").Value = 0.1822966222356014
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 10699.84005902664
$ws.Range("R14").Value = 96298.56053123974
$ws.Range("S14").Value = 0.05224507713817942
$ws.Range("T14").Value = 0.05224507713817941

$ws.Range("G15").Value = 73.43709933333334
$ws.Range("H15").Value = 220.311298
$ws.Range("I15").Value = 0.1822966222356014
$ws.Range("J15").Value = 0.1822966222356014
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 12396.16053216913
$ws.Range("R15").Value = 111565.4447895222
$ws.Range("S15").Value = 0.06052785458919721
$ws.Range("T15").Value = 0.0605278545891972

$ws.Range("G16").Value = 73.43709933333334
$ws.Range("H16").Value = 220.311298
$ws.Range("I16").Value = 0.1822966222356014
$ws.Range("J16").Value = 0.1822966222356014
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 9409.213147454031
$ws.Range("R16").Value = 84682.91832708627
$ws.Range("S16").Value = 0.04594321634589249
$ws.Range("T16").Value = 0.04594321634589248

$ws.Range("G17").Value = 73.43709933333334
$ws.Range("H17").Value = 220.311298
$ws.Range("I17").Value = 0.1822966222356014
$ws.Range("J17").Value = 0.1822966222356014
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 4829.302890790181
$ws.Range("R17").Value = 43463.72601711162
$ws.Range("S17").Value = 0.02358047416233231
$ws.Range("T17").Value = 0.0235804741623323

